# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted into the "Ajo" (garlic) sheet.
# The new record belongs between the existing row 101 and row 102 (by
# date order), so row 102 is inserted (shifting the former rows
# 102:157 down to 103:158) and populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102; this pushes the former rows
# 102-157 down to 103-158 and grows the used range to A1:R158.
$ws.Rows("102:102").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A102").Value = 5
$ws.Range("B102").Value = "Macroferia Regional de Talca"
$ws.Range("C102").Value = "Maule"
$ws.Range("D102").Value = 44452
$ws.Range("E102").Value = 7
$ws.Range("F102").Value = 100112003
$ws.Range("G102").Value = "Ajo"
$ws.Range("H102").Value = "Chino"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 300
$ws.Range("K102").Value = 15000
$ws.Range("L102").Value = 15000
$ws.Range("M102").Value = 15000
$ws.Range("N102").Value = "$/malla 10 kilos"
$ws.Range("O102").Value = "China"
$ws.Range("P102").Value = 1500
$ws.Range("Q102").Value = 10
$ws.Range("R102").Value = "Hortaliza"
